$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 224.57143
$ws.Range("I33").Value = 224.57143
$ws.Range("K33").Value = 224.57143
$ws.Range("M33").Value = 4.428570000000008
$ws.Range("H40").Value = 1221.8182
$ws.Range("J40").Value = 995
$ws.Range("L40").Value = 995
$ws.Range("N40").Value = -1345
$ws.Range("H103").Value = 537.5
$ws.Range("J103").Value = 519
$ws.Range("L103").Value = 1557
$ws.Range("N103").Value = -2729
$ws.Range("H116").Value = 5750
$ws.Range("I116").Value = 4500
$ws.Range("J116").Value = 7000
$ws.Range("K116").Value = 4500
$ws.Range("L116").Value = 7000
$ws.Range("M116").Value = -1058
$ws.Range("N116").Value = -13884
$ws.Range("H138").Value = 939.2105
$ws.Range("J138").Value = 1792
$ws.Range("L138").Value = 5376
$ws.Range("N138").Value = -15656

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3062.2
$ws.Range("I45").Value = 3062.2
$ws.Range("K45").Value = 3062.2
$ws.Range("M45").Value = -2685.2
$ws.Range("H98").Value = 19000
$ws.Range("J98").Value = 19000
$ws.Range("L98").Value = 19000
$ws.Range("N98").Value = -24990
$ws.Range("H110").Value = 9439.299999999999
$ws.Range("I110").Value = 9710.333000000001
$ws.Range("K110").Value = 9710.333000000001
$ws.Range("M110").Value = -7665.333000000001
$ws.Range("H125").Value = 74983.28999999999
$ws.Range("J125").Value = 74983.28999999999
$ws.Range("L125").Value = 74983.28999999999
$ws.Range("N125").Value = -84823.28999999999
$ws.Range("H132").Value = 2533.3333
$ws.Range("I132").Value = 2542.8572
$ws.Range("K132").Value = 7628.571599999999
$ws.Range("M132").Value = -5098.571599999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 524
$ws.Range("J22").Value = 398
$ws.Range("L22").Value = 398
$ws.Range("N22").Value = -744
$ws.Range("H99").Value = 2273.75
$ws.Range("I99").Value = 1835.2
$ws.Range("J99").Value = 3004.6667
$ws.Range("K99").Value = 1835.2
$ws.Range("L99").Value = 3004.6667
$ws.Range("M99").Value = -337.2
$ws.Range("N99").Value = -6000.6667
$ws.Range("H124").Value = 47500
$ws.Range("J124").Value = 47500
$ws.Range("L124").Value = 47500
$ws.Range("N124").Value = -57320
$ws.Range("H134").Value = 5032.8945
$ws.Range("I134").Value = 5267.467
$ws.Range("K134").Value = 15802.401
$ws.Range("M134").Value = -13267.401

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 224.33333
$ws.Range("I7").Value = 175.47058
$ws.Range("K7").Value = 175.47058
$ws.Range("M7").Value = -62.47058000000001
$ws.Range("H31").Value = 3413.7856
$ws.Range("J31").Value = 6474
$ws.Range("L31").Value = 6474
$ws.Range("N31").Value = -7064
$ws.Range("H34").Value = 3413.7856
$ws.Range("J34").Value = 6474
$ws.Range("L34").Value = 6474
$ws.Range("N34").Value = -6878
$ws.Range("H43").Value = 6330
$ws.Range("J43").Value = 6330
$ws.Range("L43").Value = 6330
$ws.Range("N43").Value = -6698
$ws.Range("H99").Value = 2976.1428
$ws.Range("I99").Value = 2976.1428
$ws.Range("K99").Value = 2976.1428
$ws.Range("M99").Value = -1478.1428
$ws.Range("H101").Value = 6330
$ws.Range("J101").Value = 6330
$ws.Range("L101").Value = 6330
$ws.Range("N101").Value = -12820
$ws.Range("H126").Value = 2976.1428
$ws.Range("I126").Value = 2976.1428
$ws.Range("K126").Value = 8928.428400000001
$ws.Range("M126").Value = -6458.428400000001
$ws.Range("H134").Value = 3875.923
$ws.Range("I134").Value = 3698.818
$ws.Range("K134").Value = 11096.454
$ws.Range("M134").Value = -8561.454000000002

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 132.4
$ws.Range("J2").Value = 171
$ws.Range("L2").Value = 1026
$ws.Range("N2").Value = -1252
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("H38").Value = 163.42857
$ws.Range("J38").Value = 54.666668
$ws.Range("L38").Value = 164.000004
$ws.Range("N38").Value = -858.000004
$ws.Range("H39").Value = 5199.933
$ws.Range("J39").Value = 5356.857
$ws.Range("L39").Value = 16070.571
$ws.Range("N39").Value = -16658.571
$ws.Range("H125").Value = 23339.666
$ws.Range("I125").Value = 23339.666
$ws.Range("K125").Value = 70018.99800000001
$ws.Range("M125").Value = -65098.99800000001
$ws.Range("H131").Value = 1742.5
$ws.Range("J131").Value = 1755.3334
$ws.Range("L131").Value = 5266.0002
$ws.Range("N131").Value = -15346.0002
$ws.Range("H132").Value = 1716.3572
$ws.Range("I132").Value = 1663
$ws.Range("J132").Value = 1849.75
$ws.Range("K132").Value = 14967
$ws.Range("L132").Value = 16647.75
$ws.Range("M132").Value = -12437
$ws.Range("N132").Value = -21707.75

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 99999.336
$ws.Range("I70").Value = 99999.5
$ws.Range("J70").Value = 99999
$ws.Range("K70").Value = 99999.5
$ws.Range("L70").Value = 99999
$ws.Range("M70").Value = -99729.5
$ws.Range("N70").Value = -100539
$ws.Range("H73").Value = 99999.336
$ws.Range("I73").Value = 99999.5
$ws.Range("J73").Value = 99999
$ws.Range("K73").Value = 99999.5
$ws.Range("L73").Value = 99999
$ws.Range("M73").Value = -99063.5
$ws.Range("N73").Value = -101871
$ws.Range("H80").Value = 2401.6667
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 2401.6667
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H122").Value = 2228.5334
$ws.Range("I122").Value = 1864.0834
$ws.Range("K122").Value = 5592.2502
$ws.Range("M122").Value = -3142.2502
$ws.Range("H126").Value = 2738.3333
$ws.Range("I126").Value = 2105.5
$ws.Range("K126").Value = 6316.5
$ws.Range("M126").Value = -3846.5
$ws.Range("H132").Value = 4702.1113
$ws.Range("I132").Value = 3184.5715
$ws.Range("K132").Value = 9553.7145
$ws.Range("M132").Value = -7023.7145

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2997
$ws.Range("I7").Value = 2997
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 2997
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -2885
$ws.Range("N7").ClearContents()
$ws.Range("H22").Value = 631.1923
$ws.Range("I22").Value = 403.73685
$ws.Range("K22").Value = 403.73685
$ws.Range("M22").Value = -108.73685
$ws.Range("H27").Value = 631.1923
$ws.Range("I27").Value = 403.73685
$ws.Range("K27").Value = 403.73685
$ws.Range("M27").Value = -296.73685
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H46").Value = 1586.6875
$ws.Range("I46").Value = 1344
$ws.Range("K46").Value = 1344
$ws.Range("M46").Value = -1156
$ws.Range("H56").Value = 8550
$ws.Range("I56").Value = 9825
$ws.Range("J56").Value = 6000
$ws.Range("K56").Value = 9825
$ws.Range("L56").Value = 6000
$ws.Range("M56").Value = -9134
$ws.Range("N56").Value = -7382
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H76").Value = 9428.666999999999
$ws.Range("J76").Value = 9428.666999999999
$ws.Range("L76").Value = 9428.666999999999
$ws.Range("N76").Value = -10104.667
$ws.Range("H79").Value = 9428.666999999999
$ws.Range("J79").Value = 9428.666999999999
$ws.Range("L79").Value = 9428.666999999999
$ws.Range("N79").Value = -11768.667
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H126").Value = 2997
$ws.Range("I126").Value = 2997
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 8991
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -6521
$ws.Range("N126").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 24000
$ws.Range("J41").Value = 24000
$ws.Range("L41").Value = 24000
$ws.Range("N41").Value = -24780
$ws.Range("H96").Value = 2000
$ws.Range("J96").Value = 2000
$ws.Range("L96").Value = 2000
$ws.Range("N96").Value = -4746
$ws.Range("H126").Value = 4415.875
$ws.Range("I126").Value = 4046.7144
$ws.Range("K126").Value = 12140.1432
$ws.Range("M126").Value = -9670.143199999999
